$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $oldText) {
        $cell.Value = $newText
    }
}

Write-Host "Done updating Recorded By column"
